# Update PLC data 2025-10-13 13:54:33
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7186
$ws.Range("C3").Value = 166168
$ws.Range("C4").Value = 157088
$ws.Range("C7").Value = 5.46
$ws.Range("C8").Value = 65.13
